$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for the new "Wins" / "Losses" / "Ties" columns (AD:AF) ---
# Copy the formatting of an existing header cell (AC1) so the new header
# cells pick up the same bold/border/centered style used by the rest of
# row 1, then fill in the header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Season record values for every data row (2-49) ---
$wins = 67
$losses = 94
$ties = 0

for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
